$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

# Row 13 "Tags" gets a new "Transcriptomics" tag inserted before "mandatory",
# which pushes the existing "mandatory" value from column D to column E.
$ws.Range("E13").Value = "mandatory"
$ws.Range("D13").Value = "Transcriptomics"

# Row 14 "Tags Term Accession Number" gets the matching NCIT accession for
# the new Transcriptomics tag.
$ws.Range("D14").Value = "http://purl.obolibrary.org/obo/NCIT_C153189"

# Row 15 "Tags Term Source REF" gets the matching term source REF (NCIT)
# for the new Transcriptomics tag.
$ws.Range("D15").Value = "NCIT"
